# Update crypto price/volume table cells per latest scrape (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.564.16'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '3.442.06'
$ws.Range("E3").Value = '  -3.82%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.10'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.65'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -7.68%  '
$ws.Range("D7").Value = '3.440.78'
$ws.Range("E7").Value = '  -3.85%  '
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.50'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.121'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -9.88%  '
$ws.Range("E12").Value = '  -7.97%  '
$ws.Range("D13").Value = '4.022.37'
$ws.Range("E13").Value = '  -3.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000181'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -11.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.52'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -9.66%  '
$ws.Range("D16").Value = '3.445.69'
$ws.Range("E16").Value = '  -3.65%  '
$ws.Range("D17").Value = '65.456.55'
$ws.Range("E17").Value = '  -1.28%  '
$ws.Range("E18").Value = '  -2.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.94'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -9.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.77'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -8.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.77'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -6.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.84'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -6.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.546'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -10.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.07'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -6.58%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '3.583.26'
$ws.Range("E26").Value = '  -3.79%  '
$ws.Range("E27").Value = '  -11.22%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.30'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -10.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.25'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -9.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.17'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -12.27%  '
$ws.Range("D32").Value = '3.446.86'
$ws.Range("E32").Value = '  -3.63%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.146'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -6.50%  '
$ws.Range("E35").Value = '  -8.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '172.23'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.35%  '
$ws.Range("E37").Value = '  -12.81%  '
$ws.Range("E38").Value = '  -11.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.53'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -7.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.80'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -13.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0776'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -8.49%  '
$ws.Range("E42").Value = '  -7.23%  '
$ws.Range("E43").Value = '  -5.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.43'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -13.92%  '
$ws.Range("E46").Value = '  -12.35%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.11'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.78'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.51'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -8.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.11'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -14.66%  '
$ws.Range("D51").Value = '2.203.08'
$ws.Range("E51").Value = '  -7.54%  '
